# Update progress figures on the "Training Dashboard" sheet to reflect
# the new as-of date of 04-Nov-2025 (previously 03-Nov-2025):
#   - PERIOD TO EXPIRE (column H) decreases by one day for rows 3-5
#   - LAST UPDATE (column I) moves from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row 3
$ws.Range("H3").Value = 342
$ws.Range("I3").Value = "'04-Nov-2025"
$ws.Range("H3").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 4
$ws.Range("H4").Value = 325
$ws.Range("I4").Value = "'04-Nov-2025"
$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 5
$ws.Range("H5").Value = 286
$ws.Range("I5").Value = "'04-Nov-2025"
$ws.Range("H5").Copy() | Out-Null
$ws.Range("I5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false
